$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1: copy the formatting (bold, border, centered) from
# the adjacent existing header cell H1, then set their text values.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$excel.CutCopyMode = $false

# Data values for I and J columns, rows 2-7
$values = @(8,9,8,9,8,7)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
